$d = $word.ActiveDocument

# Step 1: Turn paragraph 12 ("To use async await...") into the bold/italic "ASYNC AWAIT:" heading
$p12 = $d.Paragraphs(12)
$quote = [char]8217
$oldSentence = "To use async await external library should be imported because rust" + $quote + "s standard library do not support or does not come with an executor so we need to reach out to an external crate for this."
$null = $p12.Range.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "ASYNC AWAIT:", 2)
$p12 = $d.Paragraphs(12)
$p12.Range.Font.Bold = $true
$p12.Range.Font.Italic = $true

Write-Output "After step1:"
for ($i = 10; $i -le 14; $i++) {
    $p = $d.Paragraphs($i)
    Write-Output "$i : [$($p.Range.Text)]"
}

# Step 2: Insert a new paragraph after paragraph 12 containing the original sentence,
# not bold/italic, with a firstLine indent of 36pt (720 twips) and a leading tab.
$p12 = $d.Paragraphs(12)
$p12.Range.InsertParagraphAfter()
$p13 = $d.Paragraphs(13)
$insPoint = $d.Range($p13.Range.Start, $p13.Range.Start)
$insPoint.Text = [char]9 + $oldSentence
$p13 = $d.Paragraphs(13)
$p13.Range.Font.Bold = $false
$p13.Range.Font.Italic = $false
$p13.Format.FirstLineIndent = 36

Write-Output "After step2:"
for ($i = 10; $i -le 15; $i++) {
    $p = $d.Paragraphs($i)
    Write-Output "$i : [$($p.Range.Text)]  FLI=$($p.Format.FirstLineIndent) Bold=$($p.Range.Font.Bold) Italic=$($p.Range.Font.Italic)"
}
